$d = $word.ActiveDocument

# Locate the paragraph containing the "...estão mortos..." sentence. Right
# after it (in the original document) come two paragraphs that must be
# removed in their entirety:
#   1) an empty paragraph
#   2) a paragraph with the text "Os jogos de aventura foram os mais
#      vendidos na década de 1980 e início de 1990, com Myst conquistando
#      o título de jogo mais vendido de todos os tempos até The Sims
#      (Maxis Software, 2000) ultrapassá-lo em vendas."
# Both paragraphs must disappear completely (not just their text), so that
# the "mortos" paragraph is followed directly by the next (already empty)
# paragraph that precedes "A riqueza do gênero...".

$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -match [regex]::Escape('estão mortos')) {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not find the anchor paragraph containing 'estão mortos'"
}

# Delete the paragraph right after the anchor (expected to be empty).
$following = $d.Paragraphs.Item($anchorIndex + 1)
$following.Range.Delete()

# After that deletion, the paragraph with the "mais vendidos" sentence has
# shifted into the slot right after the anchor; delete it too.
$following2 = $d.Paragraphs.Item($anchorIndex + 1)
$following2.Range.Delete()
